# Update cryptos list data (Price / Volume(1h) columns), and row 12/13 coin swap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a plain-number-looking string as TEXT instead of
    # letting Excel auto-convert it to a numeric value.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "36.486.35"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.955.80"
$ws.Range("E3").Value = "  +0.76%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "243.51"
$ws.Range("E5").Value = "  +0.13%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.629"
$ws.Range("E6").Value = "  +2.96%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "60.25"
$ws.Range("E7").Value = "  +5.58%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.376"
$ws.Range("E9").Value = "  +3.83%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0784"
$ws.Range("E10").Value = "  -2.19%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.71%  "

# Row 12 - was Chainlink, now Polygon
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D12") "0.848"
$ws.Range("E12").Value = "  +5.40%  "

# Row 13 - was Polygon, now Chainlink
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "14.11"
$ws.Range("E13").Value = "  +6.25%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.242.43"
$ws.Range("E14").Value = "  +0.68%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "21.56"
$ws.Range("E15").Value = "  -0.81%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +2.01%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.952.92"
$ws.Range("E17").Value = "  +0.58%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "36.454.12"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19 - Litecoin
Set-TextValue $ws.Range("D19") "69.11"
$ws.Range("E19").Value = "  -0.19%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0854"
$ws.Range("E20").Value = "  -0.01%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "229.30"
$ws.Range("E21").Value = "  +0.81%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +2.39%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.09%  "

# Row 24 - PancakeSwap
Set-TextValue $ws.Range("D24") "2.44"
$ws.Range("E24").Value = "  +2.20%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.36"
$ws.Range("E25").Value = "  +2.60%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +6.52%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -0.65%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "160.61"
$ws.Range("E28").Value = "  +0.69%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  +19.78%  "

# Row 31 - Stellar
Set-TextValue $ws.Range("D31") "0.120"

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +3.47%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -0.66%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +6.96%  "

# Row 35 - BinanceUSD
$ws.Range("E35").Value = "  +0.04%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +2.80%  "

# Row 37 - RenderToken
Set-TextValue $ws.Range("D37") "3.37"
$ws.Range("E37").Value = "  +2.50%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -0.26%  "

# Row 39 - THORChain
Set-TextValue $ws.Range("D39") "5.43"
$ws.Range("E39").Value = "  -11.44%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +0.20%  "

# Row 42 - TrustWalletToken
$ws.Range("E42").Value = "  +1.63%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +0.37%  "

# Row 44 - InjectiveProtocol
Set-TextValue $ws.Range("D44") "15.87"
$ws.Range("E44").Value = "  +0.89%  "

# Row 45 - Maker
$ws.Range("D45").Value = "1.359.54"
$ws.Range("E45").Value = "  +1.27%  "

# Row 46 - Aave
Set-TextValue $ws.Range("D46") "88.59"
$ws.Range("E46").Value = "  +2.84%  "

# Row 47 - ARBITRUM
Set-TextValue $ws.Range("D47") "1.02"
$ws.Range("E47").Value = "  -0.05%  "

# Row 48 - FraxShare
Set-TextValue $ws.Range("D48") "7.20"
$ws.Range("E48").Value = "  +1.16%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  -0.03%  "

# Row 50 - MultiversX
Set-TextValue $ws.Range("D50") "46.05"
$ws.Range("E50").Value = "  +7.04%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.137.68"
$ws.Range("E51").Value = "  +0.84%  "
